$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 16:22"

# --- Update numeric stats for specific countries that only changed values (no reordering) ---

# Row 8: Alemania
$ws.Cells.Item(8,2).Value = 148925
$ws.Cells.Item(8,3).Value = 472
$ws.Cells.Item(8,5).Value = 44408
$ws.Cells.Item(8,7).Value = 31
$ws.Cells.Item(8,8).Value = 5117

# Row 9: Reino Unido
$ws.Cells.Item(9,2).Value = 133495
$ws.Cells.Item(9,3).Value = 4451
$ws.Cells.Item(9,5).Value = 115051
$ws.Cells.Item(9,7).Value = 763
$ws.Cells.Item(9,8).Value = 18100

# Row 57: Argentina
$ws.Cells.Item(57,4).Value = 872
$ws.Cells.Item(57,5).Value = 2120
$ws.Cells.Item(57,7).Value = 1
$ws.Cells.Item(57,8).Value = 152

# Row 68: Uzbekistan
$ws.Cells.Item(68,4).Value = 450
$ws.Cells.Item(68,5).Value = 1235

# Row 114: Sri Lanka
$ws.Cells.Item(114,2).Value = 323
$ws.Cells.Item(114,3).Value = 13
$ws.Cells.Item(114,4).Value = 105
$ws.Cells.Item(114,5).Value = 211

# Row 181: Botsuana
$ws.Cells.Item(181,2).Value = 22
$ws.Cells.Item(181,3).Value = 2
$ws.Cells.Item(181,5).Value = 21

# --- Mali moved up in the country list (now reported right after Kenia, before
#     Venezuela). Rows 120-124 shift down by one country (Venezuela, Somalia,
#     Tanzania, Vietnam each move to the next row), and row 120 now holds Mali
#     with its newly-updated figures. ---

# Row 120: now Mali (new figures)
$ws.Cells.Item(120,1).Value = "Mali"
$ws.Cells.Item(120,2).Value = 293
$ws.Cells.Item(120,3).Value = 35
$ws.Cells.Item(120,4).Value = 73
$ws.Cells.Item(120,5).Value = 203
$ws.Cells.Item(120,6).Value = 0
$ws.Cells.Item(120,7).Value = 3
$ws.Cells.Item(120,8).Value = 17

# Row 121: now Venezuela (carrying old Venezuela figures)
$ws.Cells.Item(121,1).Value = "Venezuela"
$ws.Cells.Item(121,2).Value = 288
$ws.Cells.Item(121,3).Value = 0
$ws.Cells.Item(121,4).Value = 122
$ws.Cells.Item(121,5).Value = 156
$ws.Cells.Item(121,6).Value = 4
$ws.Cells.Item(121,7).Value = 0
$ws.Cells.Item(121,8).Value = 10

# Row 122: now Somalia (carrying old Somalia figures)
$ws.Cells.Item(122,1).Value = "Somalia"
$ws.Cells.Item(122,2).Value = 286
$ws.Cells.Item(122,3).Value = 0
$ws.Cells.Item(122,4).Value = 4
$ws.Cells.Item(122,5).Value = 274
$ws.Cells.Item(122,6).Value = 2
$ws.Cells.Item(122,7).Value = 0
$ws.Cells.Item(122,8).Value = 8

# Row 123: now Tanzania (carrying old Tanzania figures)
$ws.Cells.Item(123,1).Value = "Tanzania"
$ws.Cells.Item(123,2).Value = 284
$ws.Cells.Item(123,3).Value = 30
$ws.Cells.Item(123,4).Value = 11
$ws.Cells.Item(123,5).Value = 263
$ws.Cells.Item(123,6).Value = 7
$ws.Cells.Item(123,7).Value = 0
$ws.Cells.Item(123,8).Value = 10

# Row 124: now Vietnam (carrying old Vietnam figures)
$ws.Cells.Item(124,1).Value = "Vietnam"
$ws.Cells.Item(124,2).Value = 268
$ws.Cells.Item(124,3).Value = 0
$ws.Cells.Item(124,4).Value = 223
$ws.Cells.Item(124,5).Value = 45
$ws.Cells.Item(124,6).Value = 8
$ws.Cells.Item(124,7).Value = 0
$ws.Cells.Item(124,8).Value = 0
